# ------------------------------------------------------------------
# feat: add 2022-Q1 data
# Insert a new "2022-Q1" worksheet (fund holdings detail) positioned
# between "2021-Q4" and "总计", and add a summary row for it on the
# "总计" sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" sheet as a copy of "2021-Q4" so it
#        inherits the identical header row / column styling, then
#        overwrite its data. -------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$ws = $wb.Worksheets.Item($template.Index + 1)
$ws.Name = "2022-Q1"

# Extend the two template data rows (2-3) down to 11 rows, copying
# the formatting (A column bold/centered/bordered style, plain data
# cells) of row 3 into rows 4:11.
$ws.Rows.Item(3).Copy()
$ws.Range("A4:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'270023"
$ws.Cells.Item(2,3).Value = "广发全球精选股票(QDII)"
$ws.Cells.Item(2,4).Value = "'25.53"
$ws.Cells.Item(2,5).Value = "'78.43"
$ws.Cells.Item(2,6).Value = "'3.76"
$ws.Cells.Item(2,7).Value = "'0.9599"
$ws.Cells.Item(2,8).Value = 7

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'000906"
$ws.Cells.Item(3,3).Value = "广发全球精选股票(QDII)美元现汇"
$ws.Cells.Item(3,4).Value = "'25.53"
$ws.Cells.Item(3,5).Value = "'78.43"
$ws.Cells.Item(3,6).Value = "'3.76"
$ws.Cells.Item(3,7).Value = "'0.9599"
$ws.Cells.Item(3,8).Value = 7

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'001092"
$ws.Cells.Item(4,3).Value = "广发纳斯达克生物科技指数(QDII)（人民币）"
$ws.Cells.Item(4,4).Value = "'1.34"
$ws.Cells.Item(4,5).Value = "'82.00"
$ws.Cells.Item(4,6).Value = "'3.64"
$ws.Cells.Item(4,7).Value = "'0.0488"
$ws.Cells.Item(4,8).Value = 6

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'001093"
$ws.Cells.Item(5,3).Value = "广发纳斯达克生物科技指数(QDII)（美元）"
$ws.Cells.Item(5,4).Value = "'1.34"
$ws.Cells.Item(5,5).Value = "'82.00"
$ws.Cells.Item(5,6).Value = "'3.64"
$ws.Cells.Item(5,7).Value = "'0.0488"
$ws.Cells.Item(5,8).Value = 6

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'003720"
$ws.Cells.Item(6,3).Value = "易方达标普生物科技指数（QDII-LOF）美元"
$ws.Cells.Item(6,4).Value = "'2.11"
$ws.Cells.Item(6,5).Value = "'94.00"
$ws.Cells.Item(6,6).Value = "'1.01"
$ws.Cells.Item(6,7).Value = "'0.0213"
$ws.Cells.Item(6,8).Value = 4

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'161127"
$ws.Cells.Item(7,3).Value = "易方达标普生物科技指数（QDII-LOF）人民币"
$ws.Cells.Item(7,4).Value = "'2.11"
$ws.Cells.Item(7,5).Value = "'94.00"
$ws.Cells.Item(7,6).Value = "'1.01"
$ws.Cells.Item(7,7).Value = "'0.0213"
$ws.Cells.Item(7,8).Value = 4

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'003719"
$ws.Cells.Item(8,3).Value = "易方达标普医疗保健指数(QDII-LOF) 美元"
$ws.Cells.Item(8,4).Value = "'0.51"
$ws.Cells.Item(8,5).Value = "'94.20"
$ws.Cells.Item(8,6).Value = "'1.68"
$ws.Cells.Item(8,7).Value = "'0.0086"
$ws.Cells.Item(8,8).Value = 3

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'161126"
$ws.Cells.Item(9,3).Value = "易方达标普医疗保健指数(QDII-LOF) 人民币"
$ws.Cells.Item(9,4).Value = "'0.51"
$ws.Cells.Item(9,5).Value = "'94.20"
$ws.Cells.Item(9,6).Value = "'1.68"
$ws.Cells.Item(9,7).Value = "'0.0086"
$ws.Cells.Item(9,8).Value = 3

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'096001"
$ws.Cells.Item(10,3).Value = "大成标普500等权重指数(QDII)"
$ws.Cells.Item(10,4).Value = "'3.58"
$ws.Cells.Item(10,5).Value = "'93.44"
$ws.Cells.Item(10,6).Value = "'0.22"
$ws.Cells.Item(10,7).Value = "'0.0079"
$ws.Cells.Item(10,8).Value = 8

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'013404"
$ws.Cells.Item(11,3).Value = "大成标普500等权重指数（QDII）美元"
$ws.Cells.Item(11,4).Value = "'3.58"
$ws.Cells.Item(11,5).Value = "'93.44"
$ws.Cells.Item(11,6).Value = "'0.22"
$ws.Cells.Item(11,7).Value = "'0.0079"
$ws.Cells.Item(11,8).Value = 8

# --- 2. Insert a new top data row on "总计" for 2022-Q1, pushing the
#        existing quarters down, and renumber the index column. ----
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()
$zj.Range("A3:D3").Copy()
$zj.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = "2022-Q1"
$zj.Cells.Item(2,3).Value = 10
$zj.Cells.Item(2,4).Value = 2.09

$zj.Cells.Item(3,1).Value = 1
$zj.Cells.Item(4,1).Value = 2
$zj.Cells.Item(5,1).Value = 3

# --- 3. Restore the originally-active "2020-Q4" tab as the selected
#        sheet (creating/editing the other sheets above shifts focus). ---
$wb.Worksheets.Item("2020-Q4").Activate()

